{"js": "// Rename the control/text \"-CongViecConLai\" to \"-CongViec\" in the\n// \"D\u1ef1 \u00e1n\" (Project) section of the document.\nconst body = context.document.body;\n\nconst results = body.search(\"-CongViecConLai\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const range = results.items[0];\n  range.insertText(\"-CongViec\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Rename the control/text \"-CongViecConLai\" to \"-CongViec\" in the\n# \"D\u1ef1 \u00e1n\" (Project) section of the document.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"-CongViecConLai\"\n$find.Replacement.Text = \"-CongViec\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n[void]$find.Execute([ref]\"-CongViecConLai\", [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]\"-CongViec\", [ref]2)\n"}
